# Update "部门情况202404" (department situation) sheet:
#   - row 7  (普惠业务部)   : F/G/H/I revised
#   - row 10 (普惠业务二部) : F/G/H/I revised
#   - row 11 (普惠业务三部) : brand-new row appended
#
# All data cells in this sheet store plain numbers as literal TEXT
# (OOXML t="inlineStr"), so every numeric-looking value below is written
# with the cell pre-formatted as Text ("@") — exactly what a person
# re-keying this report in real Excel would do to stop the numbers
# auto-converting — before the literal is assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).NumberFormat = "@"
    $ws.Range($range).Value = $value
}

# --- Row 7: 普惠业务部 ---------------------------------------------------
Set-TextValue "F7" "314113.65"
Set-TextValue "G7" "3134.00"
Set-TextValue "H7" "5.99"
Set-TextValue "I7" "118469.22"

# --- Row 10: 普惠业务二部 ------------------------------------------------
Set-TextValue "F10" "10209.00"
Set-TextValue "G10" "52.00"
Set-TextValue "H10" "6.27"
Set-TextValue "I10" "10209.00"

# --- Row 11 (new): 普惠业务三部 ------------------------------------------
Set-TextValue "A11" "普惠业务三部"
Set-TextValue "B11" "0.00"
Set-TextValue "C11" "0.00"
Set-TextValue "D11" "0.00"
Set-TextValue "E11" "0.00"
Set-TextValue "F11" "700.00"
Set-TextValue "G11" "4.00"
Set-TextValue "H11" "8.12"
Set-TextValue "I11" "700.00"
Set-TextValue "J11" "0.00"
Set-TextValue "K11" "0.00"
Set-TextValue "L11" "0.00"
Set-TextValue "M11" "0.00"
Set-TextValue "N11" "0.00"
Set-TextValue "O11" "0.00"
